# Daily update at 8 AM UTC
# Append the next day's row of data to the "Wins Over Time" tracker and
# move the "latest row" date formatting down to the newly appended row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previously-last row (16) was styled with the short "YYYY-MM-DD" date
# format to mark it as the latest entry. Since we're appending a new last
# row, row 16's date cell reverts to the regular datetime format used by
# all the other (non-final) rows.
$ws.Range("A16").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 17.
$ws.Range("A17").Value = 45757
$ws.Range("B17").Value = 63
$ws.Range("C17").Value = 66
$ws.Range("D17").Value = 61

# The new last row gets the short date format that marks the latest entry.
$ws.Range("A17").NumberFormat = "YYYY-MM-DD"
